# Auto-generated edit script applying the cryptos.xlsx price/volume refresh
# described in the commit "Updated cryptos list on Wed Feb 22 06:32:38 UTC 2023 with GitHub Actions".
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value into a cell as literal TEXT (never let Excel
# auto-convert digit-and-dot strings like "5.370" or "24.080.84" into a
# number, which would destroy the original text formatting/leading zeros).
# Route it through a temporary formula returning a quoted string literal,
# then flatten the formula down to its static value via copy / paste-values,
# leaving the cell a vanilla text value with no residual formula or style.
function Set-TextValue($cell, $text) {
    $escaped = $text.Replace("'", "''")
    $cell.Formula = "=""" + $text.Replace('"', '""') + """"
    $cell.Copy()
    $cell.PasteSpecial(-4163)
}

# Row 2
Set-TextValue $ws.Range("D2") "24.080.84"
$ws.Range("E2").Value = "  -3.56%  "

# Row 3
Set-TextValue $ws.Range("D3") "1.640.35"
$ws.Range("E3").Value = "  -3.60%  "

# Row 4
Set-TextValue $ws.Range("D4") "1.004"
$ws.Range("E4").Value = "  +0.30%  "

# Row 5
Set-TextValue $ws.Range("D5") "307.91"
$ws.Range("E5").Value = "  -2.61%  "

# Row 6
$ws.Range("E6").Value = "  +0.17%  "

# Row 7
Set-TextValue $ws.Range("D7") "0.3926"
$ws.Range("E7").Value = "  -1.76%  "

# Row 8
Set-TextValue $ws.Range("D8") "0.3857"
$ws.Range("E8").Value = "  -4.50%  "

# Row 9
Set-TextValue $ws.Range("D9") "1.002"
$ws.Range("E9").Value = "  -0.04%  "

# Row 10
Set-TextValue $ws.Range("D10") "1.355"
$ws.Range("E10").Value = "  -7.82%  "

# Row 11
Set-TextValue $ws.Range("D11") "48.62"
$ws.Range("E11").Value = "  -8.69%  "

# Row 12
Set-TextValue $ws.Range("D12") "0.08477"
$ws.Range("E12").Value = "  -3.78%  "

# Row 13
Set-TextValue $ws.Range("D13") "23.93"
$ws.Range("E13").Value = "  -7.92%  "

# Row 14
Set-TextValue $ws.Range("D14") "7.136"
$ws.Range("E14").Value = "  -4.64%  "

# Row 15
Set-TextValue $ws.Range("D15") "0.00001283"
$ws.Range("E15").Value = "  -5.39%  "

# Row 16
Set-TextValue $ws.Range("D16") "7.481"
$ws.Range("E16").Value = "  -6.21%  "

# Row 17
Set-TextValue $ws.Range("D17") "1.643.02"
$ws.Range("E17").Value = "  -5.60%  "

# Row 18
Set-TextValue $ws.Range("D18") "94.29"
$ws.Range("E18").Value = "  -1.61%  "

# Row 19
Set-TextValue $ws.Range("D19") "0.06943"
$ws.Range("E19").Value = "  -3.68%  "

# Row 20
Set-TextValue $ws.Range("D20") "20.84"
$ws.Range("E20").Value = "  +0.61%  "

# Row 21
Set-TextValue $ws.Range("D21") "6.939"
$ws.Range("E21").Value = "  -5.32%  "

# Row 22
$ws.Range("E22").Value = "  +0.18%  "

# Row 23
Set-TextValue $ws.Range("D23") "13.68"
$ws.Range("E23").Value = "  -4.58%  "

# Row 24
Set-TextValue $ws.Range("D24") "24.140.02"
$ws.Range("E24").Value = "  -3.29%  "

# Row 25
Set-TextValue $ws.Range("D25") "2.344"
$ws.Range("E25").Value = "  -1.67%  "

# Row 26
Set-TextValue $ws.Range("D26") "2.696"
$ws.Range("E26").Value = "  -8.88%  "

# Row 27
Set-TextValue $ws.Range("D27") "22.47"
$ws.Range("E27").Value = "  -4.73%  "

# Row 28
Set-TextValue $ws.Range("D28") "8.835"
$ws.Range("E28").Value = "  +6.00%  "

# Row 29
Set-TextValue $ws.Range("D29") "158.22"
$ws.Range("E29").Value = "  -2.79%  "

# Row 30
Set-TextValue $ws.Range("D30") "141.39"
$ws.Range("E30").Value = "  -6.24%  "

# Row 31
Set-TextValue $ws.Range("D31") "5.370"
$ws.Range("E31").Value = "  -12.50%  "

# Row 32
Set-TextValue $ws.Range("D32") "2.476"
$ws.Range("E32").Value = "  -6.05%  "

# Row 33
Set-TextValue $ws.Range("D33") "1.821.51"
$ws.Range("E33").Value = "  -6.21%  "

# Row 34
Set-TextValue $ws.Range("D34") "7.162"
$ws.Range("E34").Value = "  -0.93%  "

# Row 35
Set-TextValue $ws.Range("D35") "0.08067"
$ws.Range("E35").Value = "  -5.66%  "

# Row 36
$ws.Range("B36").Value = "VeChain"
$ws.Range("C36").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
Set-TextValue $ws.Range("D36") "0.02922"
$ws.Range("E36").Value = "  -7.57%  "

# Row 37
$ws.Range("B37").Value = "ImmutableX"
$ws.Range("C37").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
Set-TextValue $ws.Range("D37") "0.9742"
$ws.Range("E37").Value = "  -7.48%  "

# Row 38
Set-TextValue $ws.Range("D38") "0.2701"
$ws.Range("E38").Value = "  -6.25%  "

# Row 39
Set-TextValue $ws.Range("D39") "0.09224"
$ws.Range("E39").Value = "  -3.47%  "

# Row 40
Set-TextValue $ws.Range("D40") "1.466"
$ws.Range("E40").Value = "  -0.60%  "

# Row 41
Set-TextValue $ws.Range("D41") "10.03"
$ws.Range("E41").Value = "  -8.43%  "

# Row 42
Set-TextValue $ws.Range("D42") "0.7649"
$ws.Range("E42").Value = "  -7.83%  "

# Row 43
Set-TextValue $ws.Range("D43") "13.10"
$ws.Range("E43").Value = "  -6.68%  "

# Row 44
Set-TextValue $ws.Range("D44") "15.94"
$ws.Range("E44").Value = "  -7.18%  "

# Row 45
$ws.Range("B45").Value = "NEARProtocol"
$ws.Range("C45").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
Set-TextValue $ws.Range("D45") "2.483"
$ws.Range("E45").Value = "  -7.56%  "

# Row 46
$ws.Range("B46").Value = "Decentraland"
$ws.Range("C46").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
Set-TextValue $ws.Range("D46") "0.6885"
$ws.Range("E46").Value = "  -6.86%  "

# Row 47
Set-TextValue $ws.Range("D47") "4.093"
$ws.Range("E47").Value = "  -3.77%  "

# Row 48
$ws.Range("E48").Value = "  -0.04%  "

# Row 49
Set-TextValue $ws.Range("D49") "0.08398"
$ws.Range("E49").Value = "  -4.54%  "

# Row 50
Set-TextValue $ws.Range("D50") "133.99"
$ws.Range("E50").Value = "  -3.98%  "

# Row 51
Set-TextValue $ws.Range("D51") "1.260"
$ws.Range("E51").Value = "  -10.50%  "
